# Added RAD Test Cases and data for MRF.
#
# This simulates a re-run of the Katalon RAD test suite: the "Date" column
# (column B) on both result sheets is refreshed with the timestamps of the
# latest run, while every other column (Result/Execute/PaymentType/TaxType)
# keeps its existing value.

$wb = $excel.ActiveWorkbook

# --- FEINmismatch sheet (rows 2-13) ---------------------------------------
$ws1 = $wb.Worksheets.Item("FEINmismatch")

$sheet1Dates = @(
    "Mon Oct 02 16:21:14 EDT 2023",
    "Mon Oct 02 16:21:28 EDT 2023",
    "Mon Oct 02 16:21:40 EDT 2023",
    "Mon Oct 02 16:21:52 EDT 2023",
    "Mon Oct 02 16:22:04 EDT 2023",
    "Mon Oct 02 16:22:16 EDT 2023",
    "Mon Oct 02 16:22:28 EDT 2023",
    "Mon Oct 02 16:22:40 EDT 2023",
    "Mon Oct 02 16:22:52 EDT 2023",
    "Mon Oct 02 16:23:04 EDT 2023",
    "Mon Oct 02 16:23:16 EDT 2023",
    "Mon Oct 02 16:23:28 EDT 2023"
)

for ($i = 0; $i -lt $sheet1Dates.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = $sheet1Dates[$i]
}

# --- FEINSSNmismatch sheet (rows 2-17) -------------------------------------
$ws2 = $wb.Worksheets.Item("FEINSSNmismatch")

$sheet2Dates = @(
    "Mon Oct 02 16:23:40 EDT 2023",
    "Mon Oct 02 16:23:52 EDT 2023",
    "Mon Oct 02 16:24:03 EDT 2023",
    "Mon Oct 02 16:24:14 EDT 2023",
    "Mon Oct 02 16:24:26 EDT 2023",
    "Mon Oct 02 16:24:37 EDT 2023",
    "Mon Oct 02 16:24:49 EDT 2023",
    "Mon Oct 02 16:25:00 EDT 2023",
    "Mon Oct 02 16:25:11 EDT 2023",
    "Mon Oct 02 16:25:26 EDT 2023",
    "Mon Oct 02 16:25:38 EDT 2023",
    "Mon Oct 02 16:25:49 EDT 2023",
    "Mon Oct 02 16:26:01 EDT 2023",
    "Mon Oct 02 16:26:12 EDT 2023",
    "Mon Oct 02 16:26:23 EDT 2023",
    "Mon Oct 02 16:26:34 EDT 2023"
)

for ($i = 0; $i -lt $sheet2Dates.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 2).Value = $sheet2Dates[$i]
}
